$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing item quantity (Cantil): 19 -> 219
$ws.Range("B4").Value = 219

# Add new stock item in row 13: "Bandeira do Brasil" with quantity 50
$ws.Range("A13").Value = "Bandeira do Brasil"
$ws.Range("B13").Value = 50
